# Refresh the crypto price/volume snapshot (Price = column D, Volume(1h) = column E).
# Values are stored as text (inline strings in the original workbook), so each
# assignment uses a leading apostrophe to force Excel's text (quote-prefix) entry
# path instead of auto-converting number-looking text (e.g. "592.75") into a real
# number, then resets .Style back to "Normal" so no stray quote-prefix cell style
# is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = '''62.845.85'
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = '''  +2.64%  '
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = '''2.946.15'
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.Value = '''  +0.58%  '
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = '''  +0.04%  '
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = '''592.75'
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Value = '''  -0.51%  '
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = '''147.50'
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = '''  +2.41%  '
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = '''1.00'
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = '''  -0.08%  '
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = '''2.944.46'
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = '''  +0.59%  '
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = '''  +0.79%  '
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = '''7.08'
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = '''  +1.75%  '
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = '''  +5.44%  '
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = '''  +0.47%  '
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = '''  +4.55%  '
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = '''32.62'
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = '''  -1.95%  '
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = '''  -0.76%  '
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = '''3.440.75'
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = '''  +0.77%  '
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = '''62.866.57'
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = '''  +2.71%  '
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = '''6.66'
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = '''  +0.19%  '
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = '''2.948.44'
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = '''  +0.66%  '
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = '''439.70'
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = '''  +1.70%  '
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = '''13.45'
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = '''  -0.52%  '
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = '''  -0.88%  '
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = '''  -0.65%  '
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = '''11.20'
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = '''  +3.09%  '
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = '''80.76'
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = '''  -1.09%  '
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = '''  -1.81%  '
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = '''11.80'
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = '''  +0.57%  '
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = '''  -0.02%  '
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = '''  +1.54%  '
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = '''7.29'
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.Value = '''  +5.65%  '
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = '''  +0.48%  '
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = '''0.0000103'
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = '''  +17.41%  '
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = '''26.34'
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = '''  -0.89%  '
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = '''  -1.48%  '
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.Value = '''  +0.05%  '
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = '''0.992'
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.Value = '''  -1.92%  '
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = '''5.60'
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.Value = '''  -0.25%  '
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = '''3.05'
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Value = '''  +2.71%  '
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = '''49.71'
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = '''  -0.11%  '
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = '''  +1.76%  '
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = '''8.47'
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = '''  -0.61%  '
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = '''  -3.39%  '
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = '''0.278'
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = '''  -0.38%  '
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = '''39.55'
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = '''  -6.20%  '
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = '''2.706.64'
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = '''  +0.48%  '
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = '''135.56'
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = '''  +1.26%  '
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = '''0.0337'
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = '''  -1.90%  '
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = '''360.37'
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = '''  -0.60%  '
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = '''  -0.50%  '
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = '''22.69'
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = '''  -3.72%  '
$c.Style = "Normal"
